# Apply updated crypto price / volume data (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price/Coin/Link cells in this sheet are stored as text (inline strings),
# even when they look numeric (e.g. "228.07"). Force text format on the Price
# column before writing so Excel does not silently convert it to a number.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.491.52'
$ws.Range('E2').Value = '  +2.41%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.849.07'
$ws.Range('E3').Value = '  +1.98%  '

# Row 4
$ws.Range('E4').Value = '  +0.13%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.07'
$ws.Range('E5').Value = '  +0.81%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.611'
$ws.Range('E6').Value = '  +1.79%  '

# Row 7
$ws.Range('E7').Value = '  +0.12%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.42'
$ws.Range('E8').Value = '  +8.40%  '

# Row 9
$ws.Range('E9').Value = '  +5.16%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0688'
$ws.Range('E10').Value = '  +1.00%  '

# Row 11
$ws.Range('E11').Value = '  +3.46%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.118.41'
$ws.Range('E12').Value = '  +2.12%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.60'
$ws.Range('E13').Value = '  +2.74%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.847.05'
$ws.Range('E14').Value = '  +2.04%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.74'
$ws.Range('E15').Value = '  +6.54%  '

# Row 16
$ws.Range('E16').Value = '  +5.02%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '35.413.04'
$ws.Range('E17').Value = '  +2.33%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.93'
$ws.Range('E18').Value = '  +1.56%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.43'
$ws.Range('E19').Value = '  +0.24%  '

# Row 20
$ws.Range('E20').Value = '  +1.98%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.18'
$ws.Range('E21').Value = '  +7.95%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.82'
$ws.Range('E22').Value = '  +16.25%  '

# Row 24
$ws.Range('E24').Value = '  -0.75%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.21'
$ws.Range('E25').Value = '  -0.43%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.90'
$ws.Range('E26').Value = '  -0.65%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.87'
$ws.Range('E27').Value = '  +2.40%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.122'
$ws.Range('E28').Value = '  +0.94%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.486.46'
$ws.Range('E29').Value = '  +43.49%  '

# Row 30
$ws.Range('E30').Value = '  +0.18%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.33'
$ws.Range('E31').Value = '  +7.53%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.94'
$ws.Range('E32').Value = '  +3.08%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.04'
$ws.Range('E33').Value = '  +2.68%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0535'
$ws.Range('E34').Value = '  +2.02%  '

# Row 35
$ws.Range('E35').Value = '  +2.66%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.678'
$ws.Range('E36').Value = '  +3.10%  '

# Row 37
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.03'
$ws.Range('E37').Value = '  +9.51%  '

# Row 38
$ws.Range('B38').Value = 'Aave'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '88.62'
$ws.Range('E38').Value = '  +9.11%  '

# Row 39
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.08'
$ws.Range('E39').Value = '  +1.42%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.339.85'
$ws.Range('E40').Value = '  -1.98%  '

# Row 41
$ws.Range('E41').Value = '  +3.42%  '

# Row 42
$ws.Range('E42').Value = '  +2.76%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.29'
$ws.Range('E43').Value = '  +5.36%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.89'
$ws.Range('E44').Value = '  +5.25%  '

# Row 45
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.45'
$ws.Range('E45').Value = '  +0.85%  '

# Row 46
$ws.Range('B46').Value = 'MXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.82'
$ws.Range('E46').Value = '  +1.38%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0521'
$ws.Range('E47').Value = '  +3.63%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.015.90'
$ws.Range('E48').Value = '  +2.07%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.04'
$ws.Range('E49').Value = '  +4.22%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '104.43'
$ws.Range('E50').Value = '  +1.33%  '

# Row 51
$ws.Range('E51').Value = '  +0.09%  '
